$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Translate the item wording (column B) from English to German.
# Column A (item codes like "01gm") stays unchanged.
$ws.Range("B2").Value = "Die Regierung ist an Ermordungen unschuldiger Bürger und/oder bekannter Persönlichkeiten beteiligt, und hält diesen Sachverhalt geheim."
$ws.Range("B3").Value = "Die Macht von Staatsoberhäuptern ist der von kleinen anonymen Gruppen, die tatsächlich die Weltpolitik kontrollieren, untergeordnet."
$ws.Range("B4").Value = "Geheime Organisationen kommunizieren mit Außerirdischen, aber halten diese Tatsache vor der Öffentlichkeit zurück."
$ws.Range("B5").Value = "Die Verbreitung bestimmter Viren und/oder Krankheiten ist das Ergebnis der vorsätzlichen, verdeckten Aktionen einer Organisation."
$ws.Range("B6").Value = "Gruppen von Wissenschaftlern manipulieren, erfinden oder halten Beweise zurück, um die Öffentlichkeit zu täuschen."
$ws.Range("B7").Value = "Die Regierung erlaubt oder verübt selber terroristische Handlungen auf eigenem Grund und Boden und verschleiert dabei die eigene Beteiligung."
$ws.Range("B8").Value = "Ein kleiner, geheimer Personenkreis ist für das Treffen aller wichtigen Entscheidungen verantwortlich, wie z.B. in den Krieg zu ziehen."
$ws.Range("B9").Value = "Beweise für Kontakt mit Außerirdischen werden vor der Öffentlichkeit zurückgehalten."
$ws.Range("B10").Value = "Technologien, die im Stande sind Gedanken zu kontrollieren, werden an Menschen ohne deren Wissen eingesetzt."
$ws.Range("B11").Value = "Neue und fortschrittliche Technologien, die der gegenwärtigen Industrie schaden würden, werden zurückgehalten."
$ws.Range("B12").Value = "Die Regierung benutzt das Volk als Sündenbock, um die eigene Beteiligung an kriminellen Aktivitäten zu verbergen."
$ws.Range("B13").Value = "Gewisse bedeutende Ereignisse sind das Resultat der Aktivitäten einer kleinen Gruppe, die insgeheim das Weltgeschehen manipuliert."
$ws.Range("B14").Value = "Einige UFO-Sichtungen und -gerüchte werden geplant oder inszeniert, um die Öffentlichkeit von tatsächlich stattfindendem Kontakt mit Außerirdischen abzulenken."
$ws.Range("B15").Value = "Experimente, die mit neuen Medikamenten oder Technologien verbunden sind, werden regelmäßig an der Öffentlichkeit ohne deren Wissen oder Einverständnis durchgeführt."
$ws.Range("B16").Value = "Ein Großteil wichtiger Informationen wird aus Eigennutz absichtlich vor der Öffentlichkeit geheim gehalten."

# The longer German wording wraps onto more lines than the English text did,
# so several rows grow taller (matches Excel's own wrap-driven row resize).
$ws.Rows("3:3").RowHeight = 45
$ws.Rows("4:4").RowHeight = 45
$ws.Rows("6:6").RowHeight = 45
$ws.Rows("8:8").RowHeight = 45
$ws.Rows("10:10").RowHeight = 45
$ws.Rows("11:11").RowHeight = 45
$ws.Rows("12:12").RowHeight = 45
$ws.Rows("14:14").RowHeight = 45
$ws.Rows("15:15").RowHeight = 60
